# "Colocando header nos gráficos"
# Adds a header label in column A (row 1) for each data sheet, fixes
# accentuation on several row labels, and removes the per-cell bold/border
# style previously applied to the row-label cells (now only the header
# row keeps that style). Also trims the "Teto" row from the
# "Emissoes Totais" sheet and refreshes values on "Custo Total".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets 1-4 share the same layout: a "Fonte/Tecnologia" header is added
# to A1, and the row-label cells (A2:A12) lose their bold/border style.
# A few labels also get proper Portuguese accents.
# ---------------------------------------------------------------------
$sheetNames = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

$labelFixes = @{
    3  = "Gás Natural"
    4  = "Carvão"
    6  = "Óleos Comb"
    8  = "Eólica"
    11 = "Pot. Compl."
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # New header cell in A1, matching the style already used by B1:E1
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    for ($r = 2; $r -le 12; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        if ($labelFixes.ContainsKey($r)) {
            $cell.Value = $labelFixes[$r]
        }
        $cell.ClearFormats()
    }
}

# ---------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# Adds a "Período" header, fixes accents, drops the bold/border style on
# the row labels, and removes the trailing "Teto" row entirely.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

$ws5.Range("A1").Value = "Período"
$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A2").ClearFormats()

$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A3").ClearFormats()

$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
# Adds a "Tipo Expansão" header, renames the value header to "2015",
# fixes accents on the row labels, removes their bold/border style, and
# updates the cost figures.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep "2015" as text (matching the other sheets' header row), together
# with the original header style, by copying the already-text "2015"
# cell from another sheet (avoids Excel reinterpreting the digits as a
# number).
$ws1ForStyle = $wb.Worksheets.Item(1)
$ws1ForStyle.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4104)
$excel.CutCopyMode = $false

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("A2").ClearFormats()
$ws6.Range("B2").Value = 570

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("A3").ClearFormats()
$ws6.Range("B3").Value = 99
